# Re-run SGNN to annotate dialog acts following clean up work to the
# original transcripts. This updates columns I (DAMSLTag) and J (DialogAct)
# for the rows whose automatic dialog-act annotation changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column I = DAMSLTag, Column J = DialogAct
$colDAMSL = 9
$colDialogAct = 10

$updates = @(
    @{ Row = 3;   DAMSL = "ba"; DialogAct = "Appreciation" },
    @{ Row = 17;  DAMSL = "sd"; DialogAct = "Statement-non-opinion" },
    @{ Row = 19;  DAMSL = "sd"; DialogAct = "Statement-non-opinion" },
    @{ Row = 20;  DAMSL = "aa"; DialogAct = "Agree/Accept" },
    @{ Row = 24;  DAMSL = "aa"; DialogAct = "Agree/Accept" },
    @{ Row = 25;  DAMSL = "aa"; DialogAct = "Agree/Accept" },
    @{ Row = 31;  DAMSL = "sd"; DialogAct = "Statement-non-opinion" },
    @{ Row = 32;  DAMSL = "sd"; DialogAct = "Statement-non-opinion" },
    @{ Row = 39;  DAMSL = "sd"; DialogAct = "Statement-non-opinion" },
    @{ Row = 42;  DAMSL = "aa"; DialogAct = "Agree/Accept" },
    @{ Row = 47;  DAMSL = "ba"; DialogAct = "Appreciation" },
    @{ Row = 48;  DAMSL = "sv"; DialogAct = "Statement-opinion" },
    @{ Row = 63;  DAMSL = "aa"; DialogAct = "Agree/Accept" },
    @{ Row = 64;  DAMSL = "sd"; DialogAct = "Statement-non-opinion" },
    @{ Row = 70;  DAMSL = "ba"; DialogAct = "Appreciation" },
    @{ Row = 72;  DAMSL = "sv"; DialogAct = "Statement-opinion" },
    @{ Row = 82;  DAMSL = "aa"; DialogAct = "Agree/Accept" },
    @{ Row = 85;  DAMSL = "qy"; DialogAct = "Yes-No-Question" },
    @{ Row = 89;  DAMSL = "%";  DialogAct = "Uninterpretable" },
    @{ Row = 94;  DAMSL = "%";  DialogAct = "Uninterpretable" },
    @{ Row = 96;  DAMSL = "sd"; DialogAct = "Statement-non-opinion" },
    @{ Row = 104; DAMSL = "ba"; DialogAct = "Appreciation" },
    @{ Row = 105; DAMSL = "sd"; DialogAct = "Statement-non-opinion" },
    @{ Row = 107; DAMSL = "aa"; DialogAct = "Agree/Accept" },
    @{ Row = 110; DAMSL = "sd"; DialogAct = "Statement-non-opinion" },
    @{ Row = 112; DAMSL = "sd"; DialogAct = "Statement-non-opinion" },
    @{ Row = 115; DAMSL = "sv"; DialogAct = "Statement-opinion" },
    @{ Row = 127; DAMSL = "sd"; DialogAct = "Statement-non-opinion" },
    @{ Row = 136; DAMSL = "aa"; DialogAct = "Agree/Accept" },
    @{ Row = 140; DAMSL = "%";  DialogAct = "Uninterpretable" }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, $colDAMSL).Value = $u.DAMSL
    $ws.Cells.Item($u.Row, $colDialogAct).Value = $u.DialogAct
}
